$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 9
$wsALC.Range("H9").Value = 181.875
$wsALC.Range("I9").Value = 218.75
$wsALC.Range("J9").Value = 145
$wsALC.Range("K9").Value = 218.75
$wsALC.Range("L9").Value = 145
$wsALC.Range("M9").Value = -49.75
$wsALC.Range("N9").Value = -483
# ALC row 132
$wsALC.Range("H132").Value = 7250214
$wsALC.Range("I132").Value = 9528090
$wsALC.Range("K132").Value = 28584270
$wsALC.Range("M132").Value = -28581740
# ALC row 137
$wsALC.Range("H137").Value = 1465.0869
$wsALC.Range("I137").Value = 949
$wsALC.Range("K137").Value = 2847
$wsALC.Range("M137").Value = -297
# ALC row 138
$wsALC.Range("H138").Value = 1341.2527
$wsALC.Range("I138").Value = 775.32355
$wsALC.Range("J138").Value = 1656.6885
$wsALC.Range("K138").Value = 2325.97065
$wsALC.Range("L138").Value = 4970.0655
$wsALC.Range("M138").Value = 2814.02935
$wsALC.Range("N138").Value = -15250.0655
# ALC row 141
$wsALC.Range("H141").Value = 605.26666
$wsALC.Range("I141").Value = 605.26666
$wsALC.Range("J141").Value = 0
$wsALC.Range("K141").Value = 1815.79998
$wsALC.Range("L141").Value = 0
$wsALC.Range("M141").Value = 3364.20002
$wsALC.Range("N141").ClearContents()
# ARM row 74
$wsARM.Range("H74").Value = 1639.5264
$wsARM.Range("I74").Value = 781.0833
$wsARM.Range("J74").Value = 3111.1428
$wsARM.Range("K74").Value = 781.0833
$wsARM.Range("L74").Value = 3111.1428
$wsARM.Range("M74").Value = 92.91669999999999
$wsARM.Range("N74").Value = -4859.1428
# ARM row 77
$wsARM.Range("H77").Value = 1639.5264
$wsARM.Range("I77").Value = 781.0833
$wsARM.Range("J77").Value = 3111.1428
$wsARM.Range("K77").Value = 3905.4165
$wsARM.Range("L77").Value = 15555.714
$wsARM.Range("M77").Value = 462.5834999999997
$wsARM.Range("N77").Value = -24291.714
# CRP row 31
$wsCRP.Range("H31").Value = 1807.6
$wsCRP.Range("I31").Value = 1777.5
$wsCRP.Range("K31").Value = 1777.5
$wsCRP.Range("M31").Value = -1482.5
# CRP row 34
$wsCRP.Range("H34").Value = 1807.6
$wsCRP.Range("I34").Value = 1777.5
$wsCRP.Range("K34").Value = 1777.5
$wsCRP.Range("M34").Value = -1575.5
# CRP row 99
$wsCRP.Range("H99").Value = 1915.6666
$wsCRP.Range("J99").Value = 1848.4
$wsCRP.Range("L99").Value = 1848.4
$wsCRP.Range("N99").Value = -4844.4
# CRP row 126
$wsCRP.Range("H126").Value = 1915.6666
$wsCRP.Range("J126").Value = 1848.4
$wsCRP.Range("L126").Value = 5545.200000000001
$wsCRP.Range("N126").Value = -10485.2
# CUL row 122
$wsCUL.Range("H122").Value = 968
$wsCUL.Range("J122").Value = 1568
$wsCUL.Range("L122").Value = 14112
$wsCUL.Range("N122").Value = -19012
# LTW row 7
$wsLTW.Range("H7").Value = 2251
$wsLTW.Range("I7").Value = 2167.3333
$wsLTW.Range("J7").Value = 2418.3333
$wsLTW.Range("K7").Value = 2167.3333
$wsLTW.Range("L7").Value = 2418.3333
$wsLTW.Range("M7").Value = -2055.3333
$wsLTW.Range("N7").Value = -2642.3333
# LTW row 22
$wsLTW.Range("H22").Value = 787.5
$wsLTW.Range("I22").Value = 316.66666
$wsLTW.Range("J22").Value = 1070
$wsLTW.Range("K22").Value = 316.66666
$wsLTW.Range("L22").Value = 1070
$wsLTW.Range("M22").Value = -21.66665999999998
$wsLTW.Range("N22").Value = -1660
# LTW row 27
$wsLTW.Range("H27").Value = 787.5
$wsLTW.Range("I27").Value = 316.66666
$wsLTW.Range("J27").Value = 1070
$wsLTW.Range("K27").Value = 316.66666
$wsLTW.Range("L27").Value = 1070
$wsLTW.Range("M27").Value = -209.66666
$wsLTW.Range("N27").Value = -1284
# LTW row 40
$wsLTW.Range("H40").Value = 4799.364
$wsLTW.Range("I40").Value = 3098.25
$wsLTW.Range("J40").Value = 5771.4287
$wsLTW.Range("K40").Value = 3098.25
$wsLTW.Range("L40").Value = 5771.4287
$wsLTW.Range("M40").Value = -2962.25
$wsLTW.Range("N40").Value = -6043.4287
# LTW row 126
$wsLTW.Range("H126").Value = 2251
$wsLTW.Range("I126").Value = 2167.3333
$wsLTW.Range("J126").Value = 2418.3333
$wsLTW.Range("K126").Value = 6501.999899999999
$wsLTW.Range("L126").Value = 7254.999899999999
$wsLTW.Range("M126").Value = -4031.999899999999
$wsLTW.Range("N126").Value = -12194.9999
# LTW row 136
$wsLTW.Range("H136").Value = 1286.2354
$wsLTW.Range("I136").Value = 1200.871
$wsLTW.Range("J136").Value = 2168.3333
$wsLTW.Range("K136").Value = 3602.613
$wsLTW.Range("L136").Value = 6504.999899999999
$wsLTW.Range("M136").Value = -1052.613
$wsLTW.Range("N136").Value = -11604.9999
# WVR row 31
$wsWVR.Range("H31").Value = 29750
$wsWVR.Range("J31").Value = 29750
$wsWVR.Range("L31").Value = 29750
$wsWVR.Range("N31").Value = -30446
# WVR row 62
$wsWVR.Range("H62").Value = 27782308
$wsWVR.Range("I62").Value = 71434080
$wsWVR.Range("J62").Value = 3909.0908
$wsWVR.Range("K62").Value = 71434080
$wsWVR.Range("L62").Value = 3909.0908
$wsWVR.Range("M62").Value = -71433456
$wsWVR.Range("N62").Value = -5157.0908
# WVR row 65
$wsWVR.Range("H65").Value = 27782308
$wsWVR.Range("I65").Value = 71434080
$wsWVR.Range("J65").Value = 3909.0908
$wsWVR.Range("K65").Value = 357170400
$wsWVR.Range("L65").Value = 19545.454
$wsWVR.Range("M65").Value = -357167280
$wsWVR.Range("N65").Value = -25785.454
# WVR row 81
$wsWVR.Range("H81").Value = 1479.8
$wsWVR.Range("I81").Value = 799.6667
$wsWVR.Range("J81").Value = 2500
$wsWVR.Range("K81").Value = 1599.3334
$wsWVR.Range("L81").Value = 5000
$wsWVR.Range("M81").Value = -538.3334
$wsWVR.Range("N81").Value = -7122
# WVR row 84
$wsWVR.Range("H84").Value = 1479.8
$wsWVR.Range("I84").Value = 799.6667
$wsWVR.Range("J84").Value = 2500
$wsWVR.Range("K84").Value = 7996.666999999999
$wsWVR.Range("L84").Value = 25000
$wsWVR.Range("M84").Value = -2692.666999999999
$wsWVR.Range("N84").Value = -35608
# WVR row 122
$wsWVR.Range("H122").Value = 52001800
$wsWVR.Range("I122").Value = 52001800
$wsWVR.Range("K122").Value = 156005400
$wsWVR.Range("M122").Value = -156002950
# WVR row 123
$wsWVR.Range("H123").Value = 46250
$wsWVR.Range("I123").Value = 30000
$wsWVR.Range("J123").Value = 51666.668
$wsWVR.Range("K123").Value = 30000
$wsWVR.Range("L123").Value = 51666.668
$wsWVR.Range("M123").Value = -25100
$wsWVR.Range("N123").Value = -61466.668
# WVR row 126
$wsWVR.Range("H126").Value = 66668220
$wsWVR.Range("I126").Value = 90910710
$wsWVR.Range("J126").Value = 1376.25
$wsWVR.Range("K126").Value = 272732130
$wsWVR.Range("L126").Value = 4128.75
$wsWVR.Range("M126").Value = -272729660
$wsWVR.Range("N126").Value = -9068.75
# WVR row 132
$wsWVR.Range("H132").Value = 1726.625
$wsWVR.Range("I132").Value = 1435.6111
$wsWVR.Range("K132").Value = 4306.8333
$wsWVR.Range("M132").Value = -1776.8333
